# Update crypto price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "50.865.87"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "2.932.04"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "374.54"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").Value = "100.25"
$ws.Range("E6").Value = "  -3.16%  "

$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  -1.13%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  -2.47%  "

$ws.Range("D11").Value = "0.138"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "3.390.58"
$ws.Range("E13").Value = "  -1.46%  "

$ws.Range("D14").Value = "18.02"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").Value = "7.55"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.933.23"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "11.09"
$ws.Range("E17").Value = "  +49.99%  "

$ws.Range("D18").Value = "'0.990"
$ws.Range("E18").Value = "  -1.17%  "

$ws.Range("D19").Value = "50.862.24"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").Value = "3.05"
$ws.Range("E20").Value = "  -6.47%  "

$ws.Range("D21").Value = "12.41"
$ws.Range("E21").Value = "  -2.98%  "

$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("D23").Value = "68.93"
$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("D24").Value = "264.99"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").Value = "3.15"
$ws.Range("E25").Value = "  +8.52%  "

$ws.Range("D26").Value = "8.03"
$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").Value = "7.36"
$ws.Range("E27").Value = "  -2.82%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "25.57"
$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("E30").Value = "  -4.52%  "

$ws.Range("E31").Value = "  -8.17%  "

$ws.Range("D32").Value = "9.93"
$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("D33").Value = "'50.90"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").Value = "33.13"
$ws.Range("E35").Value = "  -3.78%  "

$ws.Range("D36").Value = "0.0438"
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  +4.21%  "

$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").Value = "16.42"
$ws.Range("E40").Value = "  -2.95%  "

$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  -1.59%  "

$ws.Range("D42").Value = "2.45"
$ws.Range("E42").Value = "  -4.91%  "

$ws.Range("D43").Value = "119.82"
$ws.Range("E43").Value = "  -2.19%  "

$ws.Range("D44").Value = "21.09"
$ws.Range("E44").Value = "  -2.16%  "

$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("D46").Value = "2.04"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("D47").Value = "0.268"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("D49").Value = "1.985.11"
$ws.Range("E49").Value = "  -2.28%  "

$ws.Range("D50").Value = "0.0325"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("E51").Value = "  +1.34%  "
